# Assignment 1 (Autosaved).docx - apply authoring edits via Word COM-interop
$d = $word.ActiveDocument

# Unicode helpers (avoid relying on literal smart-quote bytes through the pipe)
$ldq   = [char]0x201C   # “
$rdq   = [char]0x201D   # ”
$rsq   = [char]0x2019   # '

# ---------------------------------------------------------------------------
# 1) Heading: "Project out Line" -> "Data losses in space"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Project out Line", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Data losses in space", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) "...unmanned relay station in space. " -> "...in space combated losses. "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("to avoid losses to have unmanned relay station in space. ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "to avoid losses to have unmanned relay station in space combated losses. ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "encored and decoder data make losses significant less." -> "...decodes..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("encored and decoder data make losses significant less.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "encored and decodes data make losses significant less.", 2) | Out-Null

Write-Output "phase1 ok"

# ---------------------------------------------------------------------------
# 4) New character style "ilfuvd", based on Default Paragraph Font
# ---------------------------------------------------------------------------
$ilfuvd = $d.Styles.Add("ilfuvd", 2)
$ilfuvd.BaseStyle = $d.Styles.Item("DefaultParagraphFont")

Write-Output "style added"

# ---------------------------------------------------------------------------
# 5) Locate the paragraph ending "...significant less." to insert the new
#    content after it.
# ---------------------------------------------------------------------------
$anchorIdx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i++
    if ($p.Range.Text -like "*significant less.*") {
        $anchorIdx = $i
    }
}
Write-Output "anchorIdx=$anchorIdx"

$anchorPara = $d.Paragraphs($anchorIdx)
$anchorPara.Range.InsertParagraphAfter() | Out-Null

# ---------------------------------------------------------------------------
# 6) NASA paragraph (mixed "ilfuvd" character-styled runs)
# ---------------------------------------------------------------------------
$nasaPara = $d.Paragraphs($anchorIdx + 1)
$nasaText = " `t NASA use series of large dish antennas designed to receive the very weak from space; my proposal is compress and encored data, sent to relay station will receive signal decode the data check for losses if ok then send to next destination. I am meanly looking at compression software and Compression ratio. When we look at compress programmes like WinZip are using compression Algorithm this Quotation Wikipedia $ldq" + `
    "According to Gregory Chaitin, " + `
    "it is `"the result of putting Claude Shannon information theory and Alan Turing${rsq}s Computability theory into a cocktail shaker and shaking vigorously. The compress Algorithm is:"
$nasaPara.Range.Text = $nasaText
$nasaStart = $nasaPara.Range.Start

# "series of large dish antennas designed to receive the very weak from
#  space; my proposal is ... Wikipedia " carries the "ilfuvd" character style
$seg1 = $d.Range($nasaStart + 12, $nasaStart + 415)
$seg1.Style = $d.Styles.Item("ilfuvd")

# "According to" (offsets 415-427) is deliberately left unstyled, then
# " Gregory Chaitin, " (427-445) resumes the "ilfuvd" style
$seg2 = $d.Range($nasaStart + 427, $nasaStart + 445)
$seg2.Style = $d.Styles.Item("ilfuvd")

Write-Output "nasa paragraph populated"

# ---------------------------------------------------------------------------
# 7) "Compression ratio = Uncompressed Size/Compressed Size" paragraph
# ---------------------------------------------------------------------------
$nasaPara.Range.InsertParagraphAfter() | Out-Null
$ratioPara = $d.Paragraphs($anchorIdx + 2)
$ratioPara.Range.Text = "                                      Compression ratio = Uncompressed Size/Compressed Size"
Write-Output "ratio paragraph populated"

# ---------------------------------------------------------------------------
# 8) "Ex. You have 10mb file ..." paragraph, with the _GoBack bookmark at the
#    end (Word keeps a single _GoBack bookmark tracking the latest edit, so
#    adding it here also removes the stale one that used to sit after
#    "GitHub link ").
# ---------------------------------------------------------------------------
$ratioPara.Range.InsertParagraphAfter() | Out-Null
$exPara = $d.Paragraphs($anchorIdx + 3)
$exPara.Range.Text = "Ex. You have 10mb file after compressed 2mb ratio factor 10/2 ratio is 5:1"
$bmRange = $d.Range($exPara.Range.End - 1, $exPara.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
Write-Output "ex paragraph populated"

# ---------------------------------------------------------------------------
# 9) Trailing empty paragraph
# ---------------------------------------------------------------------------
$exPara.Range.InsertParagraphAfter() | Out-Null
Write-Output "trailing paragraph inserted"
